# test_dfs_zve.xlsx - "zve test needs to be updated as well"
#
# 1. The active selection moves from AA28 (with the view scrolled so C1 is
#    the top-left cell) to Z19 (view no longer pinned to a topLeftCell).
# 2. The child-care-cost relief term in column Z drops its 2800 EUR/child cap:
#       MIN(12*(P+N+0.96*O), 2800)   ->   12*(P+N+0.96*O)
#    for every row (single cell Z2, the shared formula group Z3:Z25, and the
#    combined-row formulas in Z22/Z23 which switch from the special
#    "two people share one cap" formula to the regular per-row formula).
# 3. Rows 22/23 (tu_id / hid) are renumbered from 31 to 32.
# All other changed cells (U/V/W/X/Z values throughout) are pure formula
# recalculation fallout from the two edits above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Selection / view -----------------------------------------------
$ws.Range("Z19").Select()

# --- 2. Drop the MIN(...,2800) cap on the Z column formulas ------------

# Z2 is its own (non-shared) formula.
$ws.Range("Z2").Formula = "=((0.6+(0.02*(T2-2005)))*(12*M2))+12*(P2+N2+0.96*O2)"

# Z3:Z21 keep the regular per-row pattern and stay one shared formula group.
$ws.Range("Z3:Z21").Formula = "=((0.6+(0.02*(T3-2005)))*(12*M3))+12*(P3+N3+0.96*O3)"

# Z22/Z23 move off the special combined-row formula onto the same regular
# per-row pattern used by the rest of the column.
$ws.Range("Z22").Formula = "=((0.6+(0.02*(T22-2005)))*(12*M22))+12*(P22+N22+0.96*O22)"
$ws.Range("Z23").Formula = "=((0.6+(0.02*(T23-2005)))*(12*M23))+12*(P23+N23+0.96*O23)"

# Z24:Z25 round out the shared group so it spans Z3:Z25 again.
$ws.Range("Z24:Z25").Formula = "=((0.6+(0.02*(T24-2005)))*(12*M24))+12*(P24+N24+0.96*O24)"

# --- 3. Renumber rows 22/23 from 31 to 32 -------------------------------
$ws.Range("A22").Value = 32
$ws.Range("B22").Value = 32
$ws.Range("A23").Value = 32
$ws.Range("B23").Value = 32
